# Applies the scheduled-runner update to the Leve profit-tracking sheets.
# Columns (per-sheet table): H currentAveragePrice, I currentAveragePriceNQ,
# J currentAveragePriceHQ, K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 10 — HQ columns reset to 0, LeveProfitHQ cleared entirely
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# Row 15
$ws.Range("H15").Value = 2396.34
$ws.Range("I15").Value = 2396.34
$ws.Range("K15").Value = 7189.02
$ws.Range("M15").Value = -7020.02

# Row 96
$ws.Range("H96").Value = 1477.3125
$ws.Range("I96").Value = 1317.5555
$ws.Range("J96").Value = 1682.7142
$ws.Range("K96").Value = 3952.6665
$ws.Range("L96").Value = 5048.142599999999
$ws.Range("M96").Value = -2579.6665
$ws.Range("N96").Value = -7794.142599999999

# Row 132
$ws.Range("H132").Value = 4424.1787
$ws.Range("I132").Value = 1319.5366
$ws.Range("J132").Value = 12910.2
$ws.Range("K132").Value = 3958.6098
$ws.Range("L132").Value = 38730.60000000001
$ws.Range("M132").Value = -1428.6098
$ws.Range("N132").Value = -43790.60000000001

# Row 137
$ws.Range("H137").Value = 2141885.8
$ws.Range("I137").Value = 2781214.8
$ws.Range("J137").Value = 1528129.8
$ws.Range("K137").Value = 8343644.399999999
$ws.Range("L137").Value = 4584389.4
$ws.Range("M137").Value = -8341094.399999999
$ws.Range("N137").Value = -4589489.4

# Row 138
$ws.Range("H138").Value = 3251.9807
$ws.Range("I138").Value = 10150
$ws.Range("J138").Value = 2677.1458
$ws.Range("K138").Value = 30450
$ws.Range("L138").Value = 8031.437399999999
$ws.Range("M138").Value = -25310
$ws.Range("N138").Value = -18311.4374

# ---------------------------------------------------------------------------
# ARM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 890609.25
$ws.Range("I32").Value = 971890.4399999999
$ws.Range("J32").Value = 16836.75
$ws.Range("K32").Value = 971890.4399999999
$ws.Range("L32").Value = 16836.75
$ws.Range("M32").Value = -971603.4399999999
$ws.Range("N32").Value = -17410.75

# ---------------------------------------------------------------------------
# BSM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 132
$ws.Range("H132").Value = 28000
$ws.Range("J132").Value = 28000
$ws.Range("L132").Value = 28000
$ws.Range("N132").Value = -38120

# ---------------------------------------------------------------------------
# CRP sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 1974.8081
$ws.Range("I31").Value = 917.1081
$ws.Range("J31").Value = 2606.016
$ws.Range("K31").Value = 917.1081
$ws.Range("L31").Value = 2606.016
$ws.Range("M31").Value = -622.1081
$ws.Range("N31").Value = -3196.016

# Row 34
$ws.Range("H34").Value = 1974.8081
$ws.Range("I34").Value = 917.1081
$ws.Range("J34").Value = 2606.016
$ws.Range("K34").Value = 917.1081
$ws.Range("L34").Value = 2606.016
$ws.Range("M34").Value = -715.1081
$ws.Range("N34").Value = -3010.016

# ---------------------------------------------------------------------------
# CUL sheet — rows 120-141 (excluding untouched row 135) lose all of their
# H:N (price/profit) cells outright; the Leve rows themselves remain.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H120:N134").ClearContents()
$ws.Range("H136:N141").ClearContents()

# ---------------------------------------------------------------------------
# GSM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 122
$ws.Range("H122").Value = 1458.2307
$ws.Range("I122").Value = 1223.5555
$ws.Range("K122").Value = 3670.6665
$ws.Range("M122").Value = -1220.6665

# ---------------------------------------------------------------------------
# WVR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 63
$ws.Range("H63").Value = 20749.666
$ws.Range("J63").Value = 20749.666
$ws.Range("L63").Value = 20749.666
$ws.Range("N63").Value = -21997.666

# Row 66
$ws.Range("H66").Value = 20749.666
$ws.Range("J66").Value = 20749.666
$ws.Range("L66").Value = 62248.99800000001
$ws.Range("N66").Value = -68488.99800000001
